# Auto-generated edit script applying scheduled price-data refresh
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1089.9584
$ws.Range("I28").Value = 1078.9286
$ws.Range("J28").Value = 1105.4
$ws.Range("K28").Value = 1078.9286
$ws.Range("L28").Value = 1105.4
$ws.Range("M28").Value = -593.9286
$ws.Range("N28").Value = -2075.4
$ws.Range("H32").Value = 3077.2144
$ws.Range("I32").Value = 2970.1428
$ws.Range("J32").Value = 3184.2856
$ws.Range("K32").Value = 2970.1428
$ws.Range("L32").Value = 3184.2856
$ws.Range("M32").Value = -2644.1428
$ws.Range("N32").Value = -3836.2856
$ws.Range("H70").Value = 2700.3333
$ws.Range("I70").Value = 2067.3333
$ws.Range("J70").Value = 3333.3333
$ws.Range("K70").Value = 6201.999899999999
$ws.Range("L70").Value = 9999.999899999999
$ws.Range("M70").Value = -5931.999899999999
$ws.Range("N70").Value = -10539.9999
$ws.Range("H73").Value = 2700.3333
$ws.Range("I73").Value = 2067.3333
$ws.Range("J73").Value = 3333.3333
$ws.Range("K73").Value = 6201.999899999999
$ws.Range("L73").Value = 9999.999899999999
$ws.Range("M73").Value = -5265.999899999999
$ws.Range("N73").Value = -11871.9999
$ws.Range("H98").Value = 4483.619
$ws.Range("I98").Value = 5645.933
$ws.Range("J98").Value = 1577.8334
$ws.Range("K98").Value = 5645.933
$ws.Range("L98").Value = 1577.8334
$ws.Range("M98").Value = -4147.933
$ws.Range("N98").Value = -4573.8334
$ws.Range("H116").Value = 803.75
$ws.Range("I116").Value = 803.75
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 803.75
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 2638.25
$ws.Range("N116").ClearContents()
$ws.Range("H122").Value = 4483.619
$ws.Range("I122").Value = 5645.933
$ws.Range("J122").Value = 1577.8334
$ws.Range("K122").Value = 16937.799
$ws.Range("L122").Value = 4733.5002
$ws.Range("M122").Value = -14487.799
$ws.Range("N122").Value = -9633.5002
$ws.Range("H138").Value = 2540.6484
$ws.Range("I138").Value = 1457.7
$ws.Range("J138").Value = 2674.3457
$ws.Range("K138").Value = 4373.1
$ws.Range("L138").Value = 8023.0371
$ws.Range("M138").Value = 766.8999999999996
$ws.Range("N138").Value = -18303.0371

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3768.12
$ws.Range("I32").Value = 3523.835
$ws.Range("J32").Value = 11666.667
$ws.Range("K32").Value = 3523.835
$ws.Range("L32").Value = 11666.667
$ws.Range("M32").Value = -3236.835
$ws.Range("N32").Value = -12240.667
$ws.Range("H52").Value = 33826.668
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 33826.668
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 33826.668
$ws.Range("N52").Value = -34462.668
$ws.Range("H61").Value = 3833136
$ws.Range("I61").Value = 5292607
$ws.Range("J61").Value = 2024.25
$ws.Range("K61").Value = 5292607
$ws.Range("L61").Value = 2024.25
$ws.Range("M61").Value = -5292395
$ws.Range("N61").Value = -2448.25
$ws.Range("H122").Value = 1532.5518
$ws.Range("I122").Value = 1381.7894
$ws.Range("J122").Value = 1819
$ws.Range("K122").Value = 4145.3682
$ws.Range("L122").Value = 5457
$ws.Range("M122").Value = -1695.3682
$ws.Range("N122").Value = -10357
$ws.Range("H136").Value = 3833136
$ws.Range("I136").Value = 5292607
$ws.Range("J136").Value = 2024.25
$ws.Range("K136").Value = 15877821
$ws.Range("L136").Value = 6072.75
$ws.Range("M136").Value = -15875271
$ws.Range("N136").Value = -11172.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 10000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 10000
$ws.Range("N29").Value = -10578
$ws.Range("H51").Value = 58500
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 58500
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 58500
$ws.Range("N51").Value = -59482
$ws.Range("H75").Value = 0
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("M75").ClearContents()
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("M78").ClearContents()
$ws.Range("N78").ClearContents()
$ws.Range("H88").Value = 9671.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 9671.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 9671.5
$ws.Range("N88").Value = -10483.5
$ws.Range("H91").Value = 9671.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 9671.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 9671.5
$ws.Range("N91").Value = -12479.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 37000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 37000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 37000
$ws.Range("N88").Value = -37812
$ws.Range("H91").Value = 37000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 37000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 37000
$ws.Range("N91").Value = -39808
$ws.Range("H107").Value = 1736720.6
$ws.Range("I107").Value = 2083815.8
$ws.Range("J107").Value = 1245
$ws.Range("K107").Value = 2083815.8
$ws.Range("L107").Value = 1245
$ws.Range("M107").Value = -2081895.8
$ws.Range("N107").Value = -5085

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 19.916666
$ws.Range("I2").Value = 16
$ws.Range("J2").Value = 27.75
$ws.Range("K2").Value = 96
$ws.Range("L2").Value = 166.5
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = -392.5
$ws.Range("H113").Value = 4242.4443
$ws.Range("I113").Value = 603.5
$ws.Range("J113").Value = 5282.143
$ws.Range("K113").Value = 1810.5
$ws.Range("L113").Value = 15846.429
$ws.Range("M113").Value = 359.5
$ws.Range("N113").Value = -20186.429
$ws.Range("H131").Value = 3632.8667
$ws.Range("I131").Value = 4428.3335
$ws.Range("J131").Value = 3434
$ws.Range("K131").Value = 13285.0005
$ws.Range("L131").Value = 10302
$ws.Range("M131").Value = -8245.000499999998
$ws.Range("N131").Value = -20382

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 23947.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 23947.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 23947.5
$ws.Range("N32").Value = -24539.5
$ws.Range("H102").Value = 3833
$ws.Range("I102").Value = 3666.182
$ws.Range("J102").Value = 4200
$ws.Range("K102").Value = 3666.182
$ws.Range("L102").Value = 4200
$ws.Range("M102").Value = -2044.182
$ws.Range("N102").Value = -7444
$ws.Range("H136").Value = 7250.9653
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 7250.9653
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 21752.8959
$ws.Range("N136").Value = -26852.8959

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6302.44
$ws.Range("I132").Value = 6415.696
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 19247.088
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -16717.088
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 18771
$ws.Range("I74").Value = 9500
$ws.Range("J74").Value = 20625.2
$ws.Range("K74").Value = 9500
$ws.Range("L74").Value = 20625.2
$ws.Range("M74").Value = -8564
$ws.Range("N74").Value = -22497.2
$ws.Range("H77").Value = 18771
$ws.Range("I77").Value = 9500
$ws.Range("J77").Value = 20625.2
$ws.Range("K77").Value = 28500
$ws.Range("L77").Value = 61875.60000000001
$ws.Range("M77").Value = -23820
$ws.Range("N77").Value = -71235.60000000001
$ws.Range("H126").Value = 2535.7646
$ws.Range("I126").Value = 1673.8667
$ws.Range("J126").Value = 9000
$ws.Range("K126").Value = 5021.6001
$ws.Range("L126").Value = 27000
$ws.Range("M126").Value = -2551.6001
$ws.Range("N126").Value = -31940
$ws.Range("H132").Value = 6647.1333
$ws.Range("I132").Value = 7139.0386
$ws.Range("J132").Value = 3449.75
$ws.Range("K132").Value = 21417.1158
$ws.Range("L132").Value = 10349.25
$ws.Range("M132").Value = -18887.1158
$ws.Range("N132").Value = -15409.25
